# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on the first data row
# (row 2, the "0c63f77c..." file) of each language sheet, reflecting a
# newer handback run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-24 09:53:02"
$zhcn.Range("H2").Value = "2016-03-24 09:53:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-24 09:53:07"
$dede.Range("H2").Value = "2016-03-24 09:53:44"
